$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37250
$ws.Range("B3").Value = 27.27952697986577
$ws.Range("B4").Value = 2.121482090215363
$ws.Range("B5").Value = 23.13
$ws.Range("B6").Value = 25.6
$ws.Range("B7").Value = 26.6
$ws.Range("B8").Value = 28.79
$ws.Range("B9").Value = 36.33
